$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44503
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 30000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 30000
$ws.Range("S3").Value = 3000

# Row 4
$ws.Range("D4").Value = 44503
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 25000
$ws.Range("S4").Value = 2500

# Row 7
$ws.Range("D7").Value = 44434
$ws.Range("M7").Value = 20

# Row 8
$ws.Range("D8").Value = 44473
$ws.Range("M8").Value = 180
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("S8").Value = 2000

# Row 9
$ws.Range("D9").Value = 44476
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 2000

# Row 10
$ws.Range("D10").Value = 44466
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 20000
$ws.Range("S10").Value = 2000

# Row 11
$ws.Range("D11").Value = 44511
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 28000
$ws.Range("O11").Value = 28000
$ws.Range("P11").Value = 28000
$ws.Range("S11").Value = 2800

# Row 12
$ws.Range("D12").Value = 44432
$ws.Range("M12").Value = 20
